$wb = $excel.ActiveWorkbook

# --- Sheet 1: quality_comparison ---
$ws1 = $wb.Worksheets.Item(1)

# Header row borders: C1 gets top+bottom thin border (no left/right),
# D1 gets top+bottom+right thin border (no left).
$c1 = $ws1.Range("C1")
$c1.Style = "Normal"
$c1.Borders.LineStyle = 1           # full thin box
$c1.Borders.Item(7).LineStyle = -4142   # clear left
$c1.Borders.Item(10).LineStyle = -4142  # clear right -> top+bottom only

$d1 = $ws1.Range("D1")
$d1.Style = "Normal"
$d1.Borders.LineStyle = 1           # full thin box
$d1.Borders.Item(7).LineStyle = -4142   # clear left -> top+bottom+right

# Anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"

# --- Sheet 2: computational_comparison ---
$ws2 = $wb.Worksheets.Item(2)

# Reuse the already-built formats from sheet 1 (copy/paste-format) so that
# the same underlying cell styles are shared instead of rebuilt from
# scratch (avoids creating duplicate/stray style entries).
$c1.Copy()
$ws2.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$ws2.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

$d1.Copy()
$ws2.Range("D1").PasteSpecial(-4122)  # xlPasteFormats
$ws2.Range("G1").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = 0

# Anonymize "fedcore" -> "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell G5
$ws2.Range("G5").ClearContents()
